# Update cryptos.xlsx Price (D) and Volume(1h) (E) columns with the latest
# scrape values, keeping the cells as plain text (matching the existing
# inline-string cells in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "307.93";        E = "0.62%" },
    @{ Row = 3;  D = "36.29";         E = "0.83%" },
    @{ Row = 4;  D = "5.048";         E = "0.70%" },
    @{ Row = 5;  D = $null;           E = "0.70%" },
    @{ Row = 6;  D = "1.968";         E = "2.43%" },
    @{ Row = 7;  D = $null;           E = "0.51%" },
    @{ Row = 8;  D = "7.869";         E = "-0.08%" },
    @{ Row = 9;  D = "0.9288";        E = "-0.23%" },
    @{ Row = 10; D = "0.1465";        E = "16.81%" },
    @{ Row = 11; D = "0.1943";        E = "1.93%" },
    @{ Row = 12; D = "0.09201";       E = "-0.41%" },
    @{ Row = 13; D = "0.03458";       E = "-1.32%" },
    @{ Row = 14; D = "0.09901";       E = "-0.14%" },
    @{ Row = 15; D = "0.001407";      E = "-0.85%" },
    @{ Row = 16; D = "0.006290";      E = "-1.79%" },
    @{ Row = 17; D = $null;           E = "6.46%" },
    @{ Row = 18; D = "3.488";         E = "6.15%" },
    @{ Row = 19; D = "0.3461";        E = "0.71%" },
    @{ Row = 20; D = "0.1293";        E = "0.73%" },
    @{ Row = 21; D = "4.835";         E = "-7.05%" },
    @{ Row = 22; D = "0.2341";        E = "-7.37%" },
    @{ Row = 23; D = "0.04400";       E = "-0.35%" },
    @{ Row = 24; D = "0.001236";      E = "0.09%" },
    @{ Row = 25; D = "0.004194";      E = "-11.31%" },
    @{ Row = 27; D = "0.0001301";     E = "-0.05%" },
    @{ Row = 39; D = "0.02047";       E = "3.98%" },
    @{ Row = 40; D = "0.05145";       E = "-1.87%" },
    @{ Row = 41; D = "0.007470";      E = "-0.98%" },
    @{ Row = 42; D = "0.01013";       E = "-0.45%" },
    @{ Row = 43; D = $null;           E = "0.20%" },
    @{ Row = 44; D = "0.002132";      E = "1.37%" },
    @{ Row = 45; D = "0.009884";      E = "-7.42%" },
    @{ Row = 46; D = "0.00006298";    E = "-1.27%" },
    @{ Row = 47; D = "0.00000000750"; E = "-0.12%" },
    @{ Row = 48; D = "64.83";         E = "-0.59%" },
    @{ Row = 49; D = "0.001600";      E = "-3.44%" },
    @{ Row = 50; D = "0.00002100";    E = "-0.12%" },
    @{ Row = 51; D = "0.0002000";     E = "-0.12%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.Value = "'" + $u.D
        $dCell.Style = "Normal"
    }

    $eCell = $ws.Cells.Item($u.Row, 5)
    $eCell.Value = "'" + $u.E
    $eCell.Style = "Normal"
}
